# "update import and export in hr and setting"
#
# The supplier sheet's "id" column (A) is dropped and all remaining columns
# shift one place to the left (B->A, C->B, ... I->H). The two trailing
# "Jan 19, 2025" date columns (J:K) are also removed. The active selection
# ends up on the (now empty) column I, matching a fresh "select whole
# column" action after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep gridlines on their (default) display state explicitly so the saved
# view doesn't pick up a stray override.
$excel.ActiveWindow.DisplayGridlines = $true

# Remove the leading "id" column -- everything to the right shifts left.
$ws.Columns("A").Delete()

# Drop the now-orphaned date columns (previously J:K, now I:J after the
# shift above) that aren't part of the exported data anymore.
$ws.Range("I1:J3").ClearContents()

# Leave the selection on column I, mirroring the post-edit state.
$ws.Columns("I").Select()
